$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated crypto price/volume/hour data (symbol list refresh).
# Each entry lists only the columns that actually changed for that row;
# D = Price, E = Volume(1h), G = Hora. Values are written with a leading
# apostrophe so Excel stores them as text (matching the original inlineStr
# cells) instead of auto-converting numeric-looking strings to numbers.
$updates = @(
    @{Row=2; D="332.20"; E="0.47%"; G="23"},
    @{Row=3; D="41.35"; E="0.60%"; G="23"},
    @{Row=4; D="5.725"; E="0.10%"; G="23"},
    @{Row=5; E="3.75%"; G="23"},
    @{Row=6; D="8.834"; E="1.09%"; G="23"},
    @{Row=7; D="4.513"; E="0.00%"; G="23"},
    @{Row=8; E="-2.41%"; G="23"},
    @{Row=9; E="-0.40%"; G="23"},
    @{Row=10; D="0.9287"; E="0.67%"; G="23"},
    @{Row=11; D="0.1252"; E="0.74%"; G="23"},
    @{Row=12; D="0.1961"; E="0.67%"; G="23"},
    @{Row=13; D="0.09348"; E="0.96%"; G="23"},
    @{Row=14; D="0.03975"; E="8.56%"; G="23"},
    @{Row=15; D="0.1065"; E="1.03%"; G="23"},
    @{Row=16; D="0.001292"; E="-0.64%"; G="23"},
    @{Row=17; D="0.006116"; E="-0.89%"; G="23"},
    @{Row=18; D="3.435"; E="1.68%"; G="23"},
    @{Row=19; G="23"},
    @{Row=20; D="9.160"; E="9.62%"; G="23"},
    @{Row=21; E="-3.80%"; G="23"},
    @{Row=22; D="0.2628"; E="-0.85%"; G="23"},
    @{Row=23; D="0.04414"; E="-0.13%"; G="23"},
    @{Row=24; D="0.001243"; E="-1.27%"; G="23"},
    @{Row=25; D="0.004380"; E="0.99%"; G="23"},
    @{Row=26; D="0.0001193"; E="-3.91%"; G="23"},
    @{Row=27; D="0.0003997"; E="0.11%"; G="23"},
    @{Row=28; G="23"},
    @{Row=29; G="23"},
    @{Row=30; G="23"},
    @{Row=31; G="23"},
    @{Row=32; G="23"},
    @{Row=33; G="23"},
    @{Row=34; G="23"},
    @{Row=35; G="23"},
    @{Row=36; G="23"},
    @{Row=37; G="23"},
    @{Row=38; G="23"},
    @{Row=39; D="0.02819"; E="0.54%"; G="23"},
    @{Row=40; D="0.05528"; E="0.57%"; G="23"},
    @{Row=41; D="0.007926"; E="3.98%"; G="23"},
    @{Row=42; D="0.1440"; E="1.10%"; G="23"},
    @{Row=43; D="0.008980"; E="-9.66%"; G="23"},
    @{Row=44; D="0.002084"; E="-6.38%"; G="23"},
    @{Row=45; D="0.01023"; E="-13.47%"; G="23"},
    @{Row=46; D="0.00007167"; E="6.56%"; G="23"},
    @{Row=47; D="0.00000000752"; E="0.27%"; G="23"},
    @{Row=48; D="0.003396"; E="15.78%"; G="23"},
    @{Row=49; D="0.002282"; E="0.21%"; G="23"},
    @{Row=50; D="0.00002104"; E="0.27%"; G="23"},
    @{Row=51; D="0.0002004"; E="0.27%"; G="23"}
)

foreach ($item in $updates) {
    $row = $item.Row
    # Use column F (Data/date, never edited) on the same row as the style
    # reference so the rewritten cells keep their original (unstyled) look.
    $refStyle = $ws.Cells.Item($row, 6).Style

    if ($item.ContainsKey("D")) {
        $cell = $ws.Cells.Item($row, 4)
        $cell.Value = "'" + $item.D
        $cell.Style = $refStyle
    }
    if ($item.ContainsKey("E")) {
        $cell = $ws.Cells.Item($row, 5)
        $cell.Value = "'" + $item.E
        $cell.Style = $refStyle
    }
    if ($item.ContainsKey("G")) {
        $cell = $ws.Cells.Item($row, 7)
        $cell.Value = "'" + $item.G
        $cell.Style = $refStyle
    }
}
